# Append the new data row (row 73) to Sheet1 and update the selection to
# match what Excel leaves behind after entering/selecting that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: A73:F73 = 800, 6, 4, 2, 2, 17
$ws.Range("A73").Value = 800
$ws.Range("B73").Value = 6
$ws.Range("C73").Value = 4
$ws.Range("D73").Value = 2
$ws.Range("E73").Value = 2
$ws.Range("F73").Value = 17

# Leave the sheet selection on the newly entered row (B73:F73), matching
# the saved workbook's cursor position.
$ws.Range("B73:F73").Select()
